$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -716
$ws.Range("N18").ClearContents()
$ws.Range("H52").Value = 387.1111
$ws.Range("I52").Value = 349.5
$ws.Range("J52").Value = 397.85715
$ws.Range("K52").Value = 1048.5
$ws.Range("L52").Value = 1193.57145
$ws.Range("M52").Value = -888.5
$ws.Range("N52").Value = -1513.57145
$ws.Range("H86").Value = 4250
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -9246
$ws.Range("H89").Value = 4250
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -46232
$ws.Range("H98").Value = 799.8
$ws.Range("I98").Value = 739.5
$ws.Range("K98").Value = 739.5
$ws.Range("M98").Value = 758.5
$ws.Range("H99").Value = 1077.1333
$ws.Range("J99").Value = 1743.875
$ws.Range("L99").Value = 5231.625
$ws.Range("N99").Value = -8227.625
$ws.Range("H122").Value = 799.8
$ws.Range("I122").Value = 739.5
$ws.Range("K122").Value = 2218.5
$ws.Range("M122").Value = 231.5
$ws.Range("H127").Value = 704.2857
$ws.Range("I127").Value = 704.2857
$ws.Range("K127").Value = 2112.8571
$ws.Range("M127").Value = 2847.1429

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H61").Value = 3999.75
$ws.Range("I61").Value = 3999.75
$ws.Range("K61").Value = 3999.75
$ws.Range("M61").Value = -3787.75
$ws.Range("H74").Value = 4527.0527
$ws.Range("I74").Value = 4534.6665
$ws.Range("J74").Value = 4498.5
$ws.Range("K74").Value = 4534.6665
$ws.Range("L74").Value = 4498.5
$ws.Range("M74").Value = -3660.6665
$ws.Range("N74").Value = -6246.5
$ws.Range("H77").Value = 4527.0527
$ws.Range("I77").Value = 4534.6665
$ws.Range("J77").Value = 4498.5
$ws.Range("K77").Value = 22673.3325
$ws.Range("L77").Value = 22492.5
$ws.Range("M77").Value = -18305.3325
$ws.Range("N77").Value = -31228.5
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 2292.8147
$ws.Range("I132").Value = 2316.24
$ws.Range("K132").Value = 6948.719999999999
$ws.Range("M132").Value = -4418.719999999999
$ws.Range("H136").Value = 3999.75
$ws.Range("I136").Value = 3999.75
$ws.Range("K136").Value = 11999.25
$ws.Range("M136").Value = -9449.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H20").Value = 3389.75
$ws.Range("I20").Value = 4779.8
$ws.Range("K20").Value = 4779.8
$ws.Range("M20").Value = -4532.8
$ws.Range("H22").Value = 381.25
$ws.Range("I22").Value = 362
$ws.Range("K22").Value = 362
$ws.Range("M22").Value = -189
$ws.Range("H86").Value = 3113.3333
$ws.Range("I86").Value = 1241.3846
$ws.Range("K86").Value = 1241.3846
$ws.Range("M86").Value = -118.3846000000001
$ws.Range("H89").Value = 3113.3333
$ws.Range("I89").Value = 1241.3846
$ws.Range("K89").Value = 6206.923000000001
$ws.Range("M89").Value = -590.9230000000007
$ws.Range("H105").Value = 2603
$ws.Range("I105").Value = 2253.75
$ws.Range("K105").Value = 2253.75
$ws.Range("M105").Value = -506.75
$ws.Range("H106").Value = 11780.2
$ws.Range("J106").Value = 11780.2
$ws.Range("L106").Value = 11780.2
$ws.Range("N106").Value = -14304.2
$ws.Range("H134").Value = 1361.3043
$ws.Range("I134").Value = 1164.1364
$ws.Range("K134").Value = 3492.4092
$ws.Range("M134").Value = -957.4092000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 973
$ws.Range("I16").Value = 973
$ws.Range("K16").Value = 973
$ws.Range("M16").Value = -686
$ws.Range("H31").Value = 5281.5186
$ws.Range("I31").Value = 2943
$ws.Range("K31").Value = 2943
$ws.Range("M31").Value = -2648
$ws.Range("H34").Value = 5281.5186
$ws.Range("I34").Value = 2943
$ws.Range("K34").Value = 2943
$ws.Range("M34").Value = -2741
$ws.Range("H113").Value = 973
$ws.Range("I113").Value = 973
$ws.Range("K113").Value = 973
$ws.Range("M113").Value = 1197
$ws.Range("H132").Value = 3793.75
$ws.Range("I132").Value = 3596.1538
$ws.Range("J132").Value = 4650
$ws.Range("K132").Value = 10788.4614
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -8258.4614
$ws.Range("N132").Value = -19010

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 163
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 181
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 543
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -889
$ws.Range("H14").Value = 502.42856
$ws.Range("I14").Value = 502.42856
$ws.Range("K14").Value = 1507.28568
$ws.Range("M14").Value = -1334.28568
$ws.Range("H103").Value = 1919.2727
$ws.Range("J103").Value = 1919.2727
$ws.Range("L103").Value = 5757.8181
$ws.Range("N103").Value = -7515.8181
$ws.Range("H132").Value = 2700
$ws.Range("I132").Value = 400
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 3600
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -1070
$ws.Range("N132").Value = -50060
$ws.Range("H139").Value = 875.3077
$ws.Range("I139").Value = 875.3077
$ws.Range("K139").Value = 2625.9231
$ws.Range("M139").Value = 2514.0769

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1821.9286
$ws.Range("I122").Value = 1773.909
$ws.Range("J122").Value = 1998
$ws.Range("K122").Value = 5321.727000000001
$ws.Range("L122").Value = 5994
$ws.Range("M122").Value = -2871.727000000001
$ws.Range("N122").Value = -10894
$ws.Range("H134").Value = 120181.336
$ws.Range("J134").Value = 120181.336
$ws.Range("L134").Value = 360544.008
$ws.Range("N134").Value = -365614.008

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3814.75
$ws.Range("I16").Value = 4222.25
$ws.Range("K16").Value = 4222.25
$ws.Range("M16").Value = -4052.25
$ws.Range("H22").Value = 1024.25
$ws.Range("I22").Value = 839
$ws.Range("J22").Value = 1333
$ws.Range("K22").Value = 839
$ws.Range("L22").Value = 1333
$ws.Range("M22").Value = -544
$ws.Range("N22").Value = -1923
$ws.Range("H27").Value = 1024.25
$ws.Range("I27").Value = 839
$ws.Range("J27").Value = 1333
$ws.Range("K27").Value = 839
$ws.Range("L27").Value = 1333
$ws.Range("M27").Value = -732
$ws.Range("N27").Value = -1547
$ws.Range("H122").Value = 3790.5
$ws.Range("I122").Value = 3783.5
$ws.Range("K122").Value = 11350.5
$ws.Range("M122").Value = -8900.5
$ws.Range("H136").Value = 3475.4167
$ws.Range("I136").Value = 2450.8333
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 7352.499899999999
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -4802.499899999999
$ws.Range("N136").Value = -18600

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 1585.7333
$ws.Range("I100").Value = 1585.7333
$ws.Range("K100").Value = 3171.4666
$ws.Range("M100").Value = -2630.4666
$ws.Range("H122").Value = 4002.2307
$ws.Range("I122").Value = 2140.6
$ws.Range("J122").Value = 5165.75
$ws.Range("K122").Value = 6421.799999999999
$ws.Range("L122").Value = 15497.25
$ws.Range("M122").Value = -3971.799999999999
$ws.Range("N122").Value = -20397.25
